# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Feb 15 19:39:13 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.833.45'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.812.85'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.33'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.551'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.52'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.94'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '3.226.36'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").Value = '2.819.40'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.893'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '51.740.11'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.36'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +7.52%  '
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.49'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '0.0₃0993'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.82'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.77'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("E25").Value = '  +4.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.66'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.28'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("B32").Value = 'VeChain'
$ws.Range("C32").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0453'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +29.10%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.72'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.24%  '
$ws.Range("E34").Value = '  +4.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0830'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.86'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.44%  '
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.11'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.61'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.56'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '126.37'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("E44").Value = '  +1.14%  '
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = '2.080.20'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.66'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("E50").Value = '  +8.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.45'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.84%  '
